$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.4
$ws.Range("G2").Value = 4.7
$ws.Range("H2").Value = 1.67
$ws.Range("I2").Value = 1.92
$ws.Range("J2").Value = 4.1
$ws.Range("K2").Value = 6.6
$ws.Range("L2").Value = 1.17
$ws.Range("N2").Value = 3.45
$ws.Range("O2").Value = 1.09
$ws.Range("P2").Value = 3.3
$ws.Range("Q2").Value = 1.28
$ws.Range("R2").Value = 2.04
$ws.Range("S2").Value = 1.68
$ws.Range("T2").Value = 1.37
$ws.Range("U2").Value = 2.8
$ws.Range("V2").Value = 2.08
$ws.Range("W2").Value = 1.27
$ws.Range("AO2").Value = 15
$ws.Range("F3").Value = 2.52
$ws.Range("G3").Value = 2.64
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 1.48
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 1.42
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.22
$ws.Range("R3").Value = 1.23
$ws.Range("S3").Value = 4.2
$ws.Range("T3").Value = 1.88
$ws.Range("U3").Value = 1.86
$ws.Range("V3").Value = 1.42
$ws.Range("W3").Value = 1.62
$ws.Range("X3").Value = 11.5
$ws.Range("Y3").Value = 12
$ws.Range("Z3").Value = 23
$ws.Range("AA3").Value = 65
$ws.Range("AC3").Value = 7.4
$ws.Range("AD3").Value = 14.5
$ws.Range("AE3").Value = 48
$ws.Range("AH3").Value = 20
$ws.Range("AI3").Value = 65
$ws.Range("AJ3").Value = 36
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 50
$ws.Range("AM3").Value = 130
$ws.Range("AN3").Value = 29
$ws.Range("AO3").Value = 50
$ws.Range("F4").Value = 1.82
$ws.Range("G4").Value = 1.87
$ws.Range("I4").Value = 6.2
$ws.Range("J4").Value = 3.45
$ws.Range("L4").Value = 1.47
$ws.Range("N4").Value = 3.2
$ws.Range("O4").Value = 1.38
$ws.Range("P4").Value = 1.74
$ws.Range("Q4").Value = 2.14
$ws.Range("R4").Value = 1.27
$ws.Range("S4").Value = 3.9
$ws.Range("T4").Value = 2.02
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.2
$ws.Range("W4").Value = 2.08
$ws.Range("AF4").Value = 11
$ws.Range("AJ4").Value = 130
$ws.Range("AN4").Value = 60
$ws.Range("F5").Value = 1.76
$ws.Range("G5").Value = 1.77
$ws.Range("H5").Value = 5.8
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 3.9
$ws.Range("L5").Value = 1.44
$ws.Range("N5").Value = 3.7
$ws.Range("O5").Value = 1.35
$ws.Range("P5").Value = 1.88
$ws.Range("Q5").Value = 2.06
$ws.Range("R5").Value = 1.33
$ws.Range("S5").Value = 3.8
$ws.Range("W5").Value = 2.3
$ws.Range("X5").Value = 13.5
$ws.Range("Y5").Value = 18
$ws.Range("AC5").Value = 8.199999999999999
$ws.Range("AI5").Value = 100
$ws.Range("AK5").Value = 19
$ws.Range("AL5").Value = 40
$ws.Range("AN5").Value = 12
$ws.Range("AO5").Value = 110
$ws.Range("F6").Value = 4.6
$ws.Range("G6").Value = 5.4
$ws.Range("H6").Value = 1.84
$ws.Range("I6").Value = 1.93
$ws.Range("J6").Value = 3.45
$ws.Range("K6").Value = 3.9
$ws.Range("L6").Value = 1.43
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 3.2
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 1.8
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.29
$ws.Range("S6").Value = 3.9
$ws.Range("T6").Value = 1.89
$ws.Range("U6").Value = 1.94
$ws.Range("V6").Value = 2.06
$ws.Range("W6").Value = 1.24
$ws.Range("Y6").Value = 9.199999999999999
$ws.Range("Z6").Value = 980
$ws.Range("AC6").Value = 9.4
$ws.Range("AJ6").Value = 140
$ws.Range("AK6").Value = 80
